$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column Y: "RF_TABNET_diff (s)" labels (rows 9-11) and "Average diff (s)" (row 12) ---
# Copy formats from the existing analogous label column (W) so borders/fonts/number
# formats match the rest of the table.
$ws.Range("W9:W12").Copy() | Out-Null
$ws.Range("Y9:Y12").PasteSpecial(-4122) | Out-Null

$ws.Range("Y9").Value = "RF_TABNET_diff (s)"
$ws.Range("Y10").Value = "RF_TABNET_diff (s)"
$ws.Range("Y11").Value = "RF_TABNET_diff (s)"

# --- Rename the "Average diff" label (W12) to "Average diff (x)" ---
$ws.Range("W12").Value = "Average diff (x)"

$ws.Range("Y12").Value = "Average diff (s)"

# --- New column Z: the "(s)" counterpart of column X (RF_TABNET_diff values) ---
$ws.Range("X9:X12").Copy() | Out-Null
$ws.Range("Z9:Z12").PasteSpecial(-4122) | Out-Null

$ws.Range("Z9").Formula = "=V9-M9"
$ws.Range("Z10").Formula = "=V10-M10"
$ws.Range("Z11").Formula = "=V11-M11"
$ws.Range("Z12").Formula = "=AVERAGE(Z9:Z11)"

$ws.Application.CutCopyMode = $false

# --- Column widths for the two new columns (best-fit like the other label/value cols) ---
$ws.Columns("Y").ColumnWidth = 17.6640625
$ws.Columns("Z").ColumnWidth = 8.6640625

# --- View: scrolled/selected further right after adding the new columns ---
$ws.Application.ActiveWindow.ScrollColumn = 20
$ws.Range("Z12").Select()
